$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.794.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.616.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.28%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.642'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.29%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -4.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.385'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.71%  '

$ws.Range("E12").Value = '  -0.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.097.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("E15").Value = '  -6.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.652.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.626.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '344.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.20%  '

$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.71%  '

$ws.Range("E25").Value = '  -2.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '572.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.44%  '

$ws.Range("E28").Value = '  -1.42%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.160'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.03%  '

$ws.Range("E36").Value = '  -2.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '151.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.43%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.31%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.67%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '155.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0594'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.69%  '

$ws.Range("E48").Value = '  +2.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.630'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.60%  '

$ws.Range("E50").Value = '  -1.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.99'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.63%  '
